$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 242 (shifts rows 242:293 down to 243:294)
$ws.Range("A242:R242").EntireRow.Insert()

# Populate the newly inserted row 242 with the new weekly data point
$ws.Cells.Item(242, 1).Value = 5
$ws.Cells.Item(242, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(242, 3).Value = "Maule"
$ws.Cells.Item(242, 4).Value = 44641
$ws.Cells.Item(242, 5).Value = 7
$ws.Cells.Item(242, 6).Value = 100112023
$ws.Cells.Item(242, 7).Value = "Brócoli"
$ws.Cells.Item(242, 8).Value = "Sin especificar"
$ws.Cells.Item(242, 9).Value = "Primera"
$ws.Cells.Item(242, 10).Value = 5000
$ws.Cells.Item(242, 11).Value = 600
$ws.Cells.Item(242, 12).Value = 600
$ws.Cells.Item(242, 13).Value = 600
$ws.Cells.Item(242, 14).Value = "$/unidad"
$ws.Cells.Item(242, 15).Value = "Región del Maule"
$ws.Cells.Item(242, 16).Value = 600
$ws.Cells.Item(242, 17).Value = 1
$ws.Cells.Item(242, 18).Value = "Hortaliza"
